$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The "bill_detail" mini-table (header + fields) moved one column to the
# right: F3:F7 -> G3:G7, and two new fields ("total", "create_at") were
# appended at G8:G9. (NB: Range.Value's getter is unreliable in this host,
# so every destination is written with its known literal text rather than
# copied from the source cell.)
# ---------------------------------------------------------------------------
$ws.Range("G3").Value = "bill_detail"
$ws.Range("G4").Value = "bill_id"
$ws.Range("G5").Value = "product_id"
$ws.Range("G6").Value = "quantity"
$ws.Range("G7").Value = "status"
$ws.Range("G8").Value = "total"
$ws.Range("G9").Value = "create_at"

$ws.Range("F3").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("F6").Value = ""
$ws.Range("F7").Value = ""

# ---------------------------------------------------------------------------
# The "product" mini-table shifted one row down: E9:E17 -> E10:E18, leaving
# E9 empty.
# ---------------------------------------------------------------------------
$ws.Range("E10").Value = "product"
$ws.Range("E11").Value = "product_id"
$ws.Range("E12").Value = "product_name"
$ws.Range("E13").Value = "manufacturer_id"
$ws.Range("E14").Value = "product_image"
$ws.Range("E15").Value = "price"
$ws.Range("E16").Value = "description"
$ws.Range("E17").Value = "product_size"
$ws.Range("E18").Value = "type_id"

$ws.Range("E9").Value = ""

# Move the active selection to G9, matching the saved workbook state.
$ws.Range("G9").Select()
